# Edit script: add "frequencies" worksheet with marker frequency data,
# update tab selections, and adjust the active tab.
# This reproduces the commit that adds inst/sample_data/simdat.xlsx's
# third sheet ("frequencies") computed from the existing genotype data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the selection on the "map" sheet (2nd sheet) to A3.
# ---------------------------------------------------------------------
$mapSheet = $wb.Worksheets.Item(2)
$mapSheet.Activate()
$mapSheet.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Add the new "frequencies" worksheet as the last sheet in the
#    workbook (after "map").
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$freqSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$freqSheet.Name = "frequencies"

# Header row
$freqSheet.Range("A1").Value = "marker"
$freqSheet.Range("B1").Value = "marker_size"
$freqSheet.Range("C1").Value = "frequency"

# Data rows (marker, marker_size, frequency)
$data = @(
    ,@("STG0001", 147, 0.015267175572519101)
    ,@("STG0001", 149, 0.0101781170483461)
    ,@("STG0001", 150, 0.017811704834605601)
    ,@("STG0001", 151, 0.015267175572519101)
    ,@("STG0010", 178, 0.015267175572519101)
    ,@("STG0010", 181, 0.017811704834605601)
    ,@("STG0010", 183, 0.015267175572519101)
    ,@("STG0010", 185, 0.0127226463104326)
    ,@("STG0016", 143, 0.017811704834605601)
    ,@("STG0016", 149, 0.0025445292620865098)
    ,@("STG0016", 152, 0.00763358778625954)
    ,@("STG0016", 155, 0.015267175572519101)
    ,@("STG0016", 172, 0.017811704834605601)
    ,@("STG0025", 216, 0.00763358778625954)
    ,@("STG0025", 217, 0.0127226463104326)
    ,@("STG0025", 221, 0.00763358778625954)
    ,@("STI0001", 196, 0.00763358778625954)
    ,@("STI0001", 199, 0.015267175572519101)
    ,@("STI0001", 205, 0.0050890585241730301)
    ,@("STI0003", 160, 0.0127226463104326)
    ,@("STI0003", 169, 0.0127226463104326)
    ,@("STI0003", 172, 0.0101781170483461)
    ,@("STI0003", 175, 0.0127226463104326)
    ,@("STI0004", 114, 0.00763358778625954)
    ,@("STI0004", 122, 0.0101781170483461)
    ,@("STI0004", 96, 0.015267175572519101)
    ,@("STI0012", 186, 0.00763358778625954)
    ,@("STI0012", 190, 0.017811704834605601)
    ,@("STI0012", 202, 0.015267175572519101)
    ,@("STI0012", 205, 0.0127226463104326)
    ,@("STI0014", 139, 0.00763358778625954)
    ,@("STI0014", 145, 0.0127226463104326)
    ,@("STI0014", 148, 0.015267175572519101)
    ,@("STI0030", 108, 0.0127226463104326)
    ,@("STI0030", 111, 0.0101781170483461)
    ,@("STI0030", 114, 0.00763358778625954)
    ,@("STI0030", 117, 0.015267175572519101)
    ,@("STI0030", 120, 0.0050890585241730301)
    ,@("STI0032", 137, 0.0127226463104326)
    ,@("STI0032", 143, 0.015267175572519101)
    ,@("STI0033", 131, 0.017811704834605601)
    ,@("STM0031", 186, 0.0127226463104326)
    ,@("STM0031", 204, 0.0101781170483461)
    ,@("STM0037", 102, 0.017811704834605601)
    ,@("STM0037", 107, 0.0127226463104326)
    ,@("STM0037", 90, 0.00763358778625954)
    ,@("STM0037", 92, 0.0203562340966921)
    ,@("STM0037", 94, 0.0127226463104326)
    ,@("STM1052", 227, 0.015267175572519101)
    ,@("STM1052", 236, 0.0127226463104326)
    ,@("STM1052", 245, 0.0101781170483461)
    ,@("STM1053", 181, 0.015267175572519101)
    ,@("STM1053", 187, 0.00763358778625954)
    ,@("STM1053", 190, 0.017811704834605601)
    ,@("STM1053", 191, 0.00763358778625954)
    ,@("STM1064", 207, 0.00763358778625954)
    ,@("STM1064", 209, 0.0127226463104326)
    ,@("STM1064", 214, 0.0127226463104326)
    ,@("STM1104", 185, 0.00763358778625954)
    ,@("STM1104", 188, 0.0127226463104326)
    ,@("STM1104", 189, 0.015267175572519101)
    ,@("STM1104", 192, 0.015267175572519101)
    ,@("STM1106", 170, 0.0101781170483461)
    ,@("STM1106", 173, 0.0127226463104326)
    ,@("STM1106", 176, 0.0101781170483461)
    ,@("STM1106", 179, 0.015267175572519101)
    ,@("STM5114", 305, 0.0101781170483461)
    ,@("STM5114", 308, 0.0127226463104326)
    ,@("STM5114", 314, 0.0127226463104326)
    ,@("STM5121", 300, 0.022900763358778602)
    ,@("STM5121", 307, 0.0127226463104326)
    ,@("STM5127", 259, 0.0101781170483461)
    ,@("STM5127", 264, 0.0127226463104326)
    ,@("STM5127", 270, 0.00763358778625954)
    ,@("STM5127", 277, 0.015267175572519101)
    ,@("STM5127", 286, 0.015267175572519101)
    ,@("STM5127", 289, 0.0101781170483461)
    ,@("STPoAc58", 246, 0.0127226463104326)
    ,@("STPoAc58", 248, 0.0101781170483461)
    ,@("STPoAc58", 250, 0.0101781170483461)
    ,@("STPoAc58", 256, 0.00763358778625954)
)

$r = 2
foreach ($row in $data) {
    $freqSheet.Cells.Item($r, 1).Value = $row[0]
    $freqSheet.Cells.Item($r, 2).Value = $row[1]
    $freqSheet.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Column widths (best-effort match to the original workbook's precise
# pixel-derived widths; the COM layer quantizes to character-width
# sixths so we pick the input that rounds to the closest achievable
# value).
$freqSheet.Columns.Item(1).ColumnWidth = 11.8
$freqSheet.Columns.Item(2).ColumnWidth = 12.7
$freqSheet.Columns.Item(3).ColumnWidth = 11.5

# Page setup to match the other worksheets (portrait orientation).
$freqSheet.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3. Make "frequencies" the active sheet/tab with A3 selected, which
#    also clears tabSelected from the previously active sheet.
# ---------------------------------------------------------------------
$freqSheet.Activate()
$freqSheet.Range("A3").Select() | Out-Null
